# Insert 4 new rows at the top of the "Piña Caramelo" data block (row 893)
# and populate them with a new weekly entry (Fecha 2022-01-17 / serial 44578),
# pushing every subsequent row down by 4 (the rest of the sheet/data is
# unaffected in content, only shifted in row position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows above row 893 (one at a time so every existing
# row - including the last ones - shifts down correctly).
$ws.Rows.Item(893).Insert()
$ws.Rows.Item(893).Insert()
$ws.Rows.Item(893).Insert()
$ws.Rows.Item(893).Insert()

# Common values shared by the four new rows.
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$fecha = 44578
$codreg = 13
$tipo = "Fruta"
$productoId = 100108
$producto = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria = "Piña"
$variedad = "Caramelo"
$origen = "Ecuador"

# Row 893 - Especial
$ws.Range("A893").Value() = $mercadoId
$ws.Range("B893").Value() = $mercado
$ws.Range("C893").Value() = $region
$ws.Range("D893").Value() = $fecha
$ws.Range("E893").Value() = $codreg
$ws.Range("F893").Value() = $tipo
$ws.Range("G893").Value() = $productoId
$ws.Range("H893").Value() = $producto
$ws.Range("I893").Value() = $categoriaId
$ws.Range("J893").Value() = $categoria
$ws.Range("K893").Value() = $variedad
$ws.Range("L893").Value() = "Especial"
$ws.Range("M893").Value() = 25
$ws.Range("N893").Value() = 14000
$ws.Range("O893").Value() = 15000
$ws.Range("P893").Value() = 14600
$ws.Range("Q893").Value() = '$/caja 10 unidades'
$ws.Range("R893").Value() = $origen
$ws.Range("S893").Value() = 1460
$ws.Range("T893").Value() = 10

# Row 894 - Primera
$ws.Range("A894").Value() = $mercadoId
$ws.Range("B894").Value() = $mercado
$ws.Range("C894").Value() = $region
$ws.Range("D894").Value() = $fecha
$ws.Range("E894").Value() = $codreg
$ws.Range("F894").Value() = $tipo
$ws.Range("G894").Value() = $productoId
$ws.Range("H894").Value() = $producto
$ws.Range("I894").Value() = $categoriaId
$ws.Range("J894").Value() = $categoria
$ws.Range("K894").Value() = $variedad
$ws.Range("L894").Value() = "Primera"
$ws.Range("M894").Value() = 30
$ws.Range("N894").Value() = 14000
$ws.Range("O894").Value() = 15000
$ws.Range("P894").Value() = 14667
$ws.Range("Q894").Value() = '$/caja 12 unidades'
$ws.Range("R894").Value() = $origen
$ws.Range("S894").Value() = 1222
$ws.Range("T894").Value() = 12

# Row 895 - Segunda
$ws.Range("A895").Value() = $mercadoId
$ws.Range("B895").Value() = $mercado
$ws.Range("C895").Value() = $region
$ws.Range("D895").Value() = $fecha
$ws.Range("E895").Value() = $codreg
$ws.Range("F895").Value() = $tipo
$ws.Range("G895").Value() = $productoId
$ws.Range("H895").Value() = $producto
$ws.Range("I895").Value() = $categoriaId
$ws.Range("J895").Value() = $categoria
$ws.Range("K895").Value() = $variedad
$ws.Range("L895").Value() = "Segunda"
$ws.Range("M895").Value() = 30
$ws.Range("N895").Value() = 14000
$ws.Range("O895").Value() = 15000
$ws.Range("P895").Value() = 14500
$ws.Range("Q895").Value() = '$/caja 14 unidades'
$ws.Range("R895").Value() = $origen
$ws.Range("S895").Value() = 1036
$ws.Range("T895").Value() = 14

# Row 896 - Tercera
$ws.Range("A896").Value() = $mercadoId
$ws.Range("B896").Value() = $mercado
$ws.Range("C896").Value() = $region
$ws.Range("D896").Value() = $fecha
$ws.Range("E896").Value() = $codreg
$ws.Range("F896").Value() = $tipo
$ws.Range("G896").Value() = $productoId
$ws.Range("H896").Value() = $producto
$ws.Range("I896").Value() = $categoriaId
$ws.Range("J896").Value() = $categoria
$ws.Range("K896").Value() = $variedad
$ws.Range("L896").Value() = "Tercera"
$ws.Range("M896").Value() = 30
$ws.Range("N896").Value() = 14000
$ws.Range("O896").Value() = 15000
$ws.Range("P896").Value() = 14333
$ws.Range("Q896").Value() = '$/caja 16 unidades'
$ws.Range("R896").Value() = $origen
$ws.Range("S896").Value() = 896
$ws.Range("T896").Value() = 16

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
